# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that hold tables) get switched from the
#    deck's custom "Table_0" style to the built-in table style
#    {1D0FD5D6-4EC4-49D2-8CAF-63CA355F7BE4}. Table styles can't be assigned
#    through the Style property directly (PowerPoint throws "Table styles
#    cannot be assigned through a property"), so Table.ApplyStyle(...) is
#    used instead.
#
# 2) The deck's theme colour palette is switched from the "Integral" /
#    "Red Violet" scheme to the default "Office" scheme. The 12 theme
#    colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) are exposed
#    as SlideMaster.ColorScheme.Colors(1..12).

$p = $ppt.ActivePresentation

# --- 1) Re-style every table in the deck -----------------------------
$newTableStyleId = "{1D0FD5D6-4EC4-49D2-8CAF-63CA355F7BE4}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Swap the theme colour scheme back to the default "Office" one -
function Set-RGBFromHex($colorScheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$masterColorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    Set-RGBFromHex $masterColorScheme $i $officeThemeColors[$i - 1]
}
